# 2nd May Data Refresh - master-reg_center_user
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_user")

# regcntr_id corrections from the refreshed data pull
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Reposition the view the way the author left it after editing (scrolled
# down to row 13, with C19 as the active/selected cell)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()
